# Ajout EDT S6 25-26.
# Fill in the room/location column (F) for each course occurrence, fix up
# the last week's session (date + day moved from Wed to Fri, and the
# course type switched from "cours" to "controle"), and remove the extra
# trailing week (old rows 25-27) that is no longer part of the schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Room assignments (column F) for existing sessions ---
$ws.Range("F3").Value  = "U3-Amphi"
$ws.Range("F4").Value  = "U3-Amphi"
$ws.Range("F7").Value  = "U3-Amphi"
$ws.Range("F8").Value  = "U3-Amphi"
$ws.Range("F11").Value = "U3-Amphi"
$ws.Range("F12").Value = "U3-Amphi"
$ws.Range("F15").Value = "U3-110"
$ws.Range("F16").Value = "U3-110"
$ws.Range("F19").Value = "U3-4"
$ws.Range("F20").Value = "U3-Amphi"

# --- Week of row 22 moves from Wednesday 26 Nov to Friday 28 Nov ---
$ws.Range("A22").Value = 45989
$ws.Range("B22").Value = "vendredi"

# --- Rows 23/24 become the "controle" (exam) session instead of "cours" ---
$ws.Range("A23").Value = "FSQTEL (KRTX9AB1)"
$ws.Range("F23").Value = "U3-Amphi"
$ws.Range("H23").Value = "TYPE_CONTROLE"

$ws.Range("A24").Value = "FSQTEL (KRTX9AB1)"
$ws.Range("F24").Value = "U3-Amphi"
$ws.Range("H24").Value = "TYPE_CONTROLE"

# --- Remove the now-unneeded trailing week (former rows 25-27) ---
$ws.Rows("25:27").Delete()
